# Apply the "Retanqueo multiple analisis de credito CCS" data edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RetanqueoMultiple")
$ws.Activate()

# Row 2 data updates (order matches shared-string allocation order in the target file)
$ws.Range("Q2").Value = '"JHON FREDY"'
$ws.Range("A2").Value = '"10002426"'
$ws.Range("C2").Value = '"P.A COLPENSIONES"'
$ws.Range("O2").Value = '"830000"'
$ws.Range("W2").Value = '"100000"'
$ws.Range("X2").Value = '"730000"'

# Update the view: active selection (matches the saved sheetView in the target)
$ws.Range("O10").Select()
